# Update "想去人数" (want-to-go count) figures that changed between the
# previous data pull and the regenerated one (456a3b4).
#
# Sheet "展览" (Exhibition):
#   F3: 2861 -> 2870
#   F5: 17   -> 40
#
# Sheet "全部类型" (All types) contains the same two events duplicated:
#   F7:  2861 -> 2870
#   F10: 17   -> 40

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 2870
$wsExhibit.Range("F5").Value = 40

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F7").Value = 2870
$wsAll.Range("F10").Value = 40
